$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few of the refreshed price strings are numeric-looking with a trailing
# zero (e.g. "301.00"). Excel.Range.Value auto-detects such text as a number
# and would silently drop the trailing zero, so pin those specific cells to
# the Text format before assigning them, matching how the source data is
# stored (plain text, not a number).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = '43.146.05'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '2.307.05'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '301.00'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").Value = '98.05'
$ws.Range("E6").Value = '  -2.21%  '
$ws.Range("D7").Value = '0.519'
$ws.Range("E7").Value = '  +2.80%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '35.89'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").Value = '0.0793'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("E13").Value = '  -4.10%  '
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '2.664.61'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '2.274.70'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  -1.67%  '
$ws.Range("D18").Value = '43.018.03'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  +2.11%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("D22").Value = '68.34'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").Value = '238.10'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '0.991'
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("E26").Value = '  -1.37%  '
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").Value = '25.20'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '166.39'
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").Value = '9.16'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  -12.60%  '
$ws.Range("D32").Value = '33.15'
$ws.Range("E32").Value = '  -4.55%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("D35").Value = '18.21'
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").Value = '0.0691'
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").Value = '1.79'
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("D42").Value = '2.76'
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("D43").Value = '2.011.69'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").Value = '2.20'
$ws.Range("E45").Value = '  -6.36%  '
$ws.Range("D46").Value = '10.30'
$ws.Range("D47").Value = '17.54'
$ws.Range("E47").Value = '  -0.34%  '
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("D49").Value = '54.41'
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").Value = '2.536.53'
$ws.Range("E51").Value = '  -1.02%  '
